$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = $true
$ws.Range("F3").Value = $true
$ws.Range("F4").Value = $false
$ws.Range("F5").Value = $true
